# Adding a new "# of Diffs" column at the front of the report (both the
# normal layout in rows 1-10 and the interlaced layout in rows 11-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts the existing
# columns A:G -> B:H, preserving their values/number formats/widths.
$ws.Columns.Item(1).Insert()

# Header cells for the new column, in both table layouts. Copy the
# neighbouring header's formatting (bold style) onto the new header cell.
$ws.Range("A1").Value = "# of Diffs"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

$ws.Range("A11").Value = "# of Diffs"
$ws.Range("B11").Copy()
$ws.Range("A11").PasteSpecial(-4122)

# Body cells: every data row gets a diff-count of 0.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = 0
}
for ($r = 12; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = 0
}

# Match the new column's width to the rest of the (bestFit) columns.
$ws.Columns.Item(1).ColumnWidth = 13.16666666666667

# Re-point the filter database defined name at the widened range.
$name = $ws.Names.Item(1)
$name.RefersTo = "=Sheet1!`$A`$1:`$H`$10"

# Re-apply the AutoFilter across the widened range (A1:H20 instead of A1:G20).
$ws.AutoFilterMode = $false
$ws.Range("A1:H20").AutoFilter()

"edit complete"
